# Update Baghdati disability_prevalence.xlsx:
#  - retitle the header row (now merged A1:I1) with the new description
#  - insert a new data row ("family with disabilities Persons") above the
#    existing data row (retitled "disabilities Persons") with refreshed
#    figures for 2017-2024
#  - re-flow the Source row down one slot (now merged A6:H6)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlThin = 2
$xlCenter = -4108
$xlLeft = -4131
$xlNone = -4142

# ---------------------------------------------------------------------
# 1. Make room: insert a fresh row above the old data row (row 4).
#    This pushes the old "Number of disability persons" row to row 5 and
#    the old "Source:" row to row 6, carrying formatting with them.
# ---------------------------------------------------------------------
$ws.Rows("4").Insert()

# ---------------------------------------------------------------------
# 2. Row 1 - title (now merged A1:I1)
# ---------------------------------------------------------------------
$r1 = $ws.Range("A1:I1")
$r1.Borders.LineStyle = $xlNone
$r1.Value2 = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Baghdati Municipality"
$r1.Merge()
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.Bold = $true
$r1.Font.Underline = $false
$r1.HorizontalAlignment = $xlCenter
$r1.VerticalAlignment = $xlCenter
$r1.WrapText = $true
$r1.Interior.Pattern = $xlNone
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------
# 3. Row 2 - "(End of year, persons)" (unchanged text/format, default height)
# ---------------------------------------------------------------------
$a2 = $ws.Range("A2")
$a2.Font.Name = "Arial"
$a2.Font.Size = 10
$a2.Font.Bold = $false
$a2.Font.Underline = $false
$a2.Interior.ThemeColor = 0
$a2.Interior.Pattern = 1
$a2.Borders.LineStyle = $xlNone
$ws.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------
# 4. Row 3 - blank A3 + year headers: untouched, formatting already correct
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5. Row 4 (NEW) - "family with disabilities Persons"
# ---------------------------------------------------------------------
$a4 = $ws.Range("A4")
$a4.Value2 = "family with disabilities Persons "
$row4vals = @(753, 717, 622, 629, 630, 624, 635, 635)
for ($i = 0; $i -lt $row4vals.Length; $i++) {
    $ws.Cells.Item(4, 2 + $i).Value2 = $row4vals[$i]
}
$ws.Rows.Item(4).RowHeight = 24.75

$a4.Font.Name = "Arial"
$a4.Font.Size = 10
$a4.Font.Bold = $false
$a4.Font.Underline = $false
$a4.Interior.ThemeColor = 0
$a4.Interior.Pattern = 1
$a4.HorizontalAlignment = $xlLeft
$a4.VerticalAlignment = $xlCenter
$a4.WrapText = $true
$a4.Borders.LineStyle = $xlNone
$a4.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$a4.Borders.Item($xlEdgeTop).Weight = $xlThin
$a4.Borders.Item($xlEdgeTop).ColorIndex = -4105

$nums4 = $ws.Range("B4:I4")
$nums4.NumberFormat = "#\ ##0"
$nums4.Font.Name = "Arial"
$nums4.Font.Size = 10
$nums4.Font.Bold = $false
$nums4.Interior.ThemeColor = 0
$nums4.Interior.Pattern = 1
$nums4.HorizontalAlignment = $xlLeft
$nums4.WrapText = $false
$nums4.Borders.LineStyle = $xlNone

# ---------------------------------------------------------------------
# 6. Row 5 (was row 4) - "disabilities Persons"
# ---------------------------------------------------------------------
$a5 = $ws.Range("A5")
$a5.Value2 = "disabilities Persons "
$row5vals = @(852, 809, 708, 715, 712, 706, 720, 724)
for ($i = 0; $i -lt $row5vals.Length; $i++) {
    $ws.Cells.Item(5, 2 + $i).Value2 = $row5vals[$i]
}
$ws.Rows.Item(5).RowHeight = 21

$a5.Font.Name = "Arial"
$a5.Font.Size = 10
$a5.Font.Bold = $false
$a5.Font.Underline = $false
$a5.Interior.ThemeColor = 0
$a5.Interior.Pattern = 1
$a5.HorizontalAlignment = $xlLeft
$a5.VerticalAlignment = $xlCenter
$a5.WrapText = $true
$a5.Borders.LineStyle = $xlNone
$a5.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$a5.Borders.Item($xlEdgeBottom).Weight = $xlThin
$a5.Borders.Item($xlEdgeBottom).ColorIndex = -4105

$nums5 = $ws.Range("B5:I5")
$nums5.NumberFormat = "#\ ##0"
$nums5.Font.Name = "Arial"
$nums5.Font.Size = 10
$nums5.Font.Bold = $false
$nums5.Interior.ThemeColor = 0
$nums5.Interior.Pattern = 1
$nums5.HorizontalAlignment = $xlLeft
$nums5.WrapText = $false
$nums5.Borders.LineStyle = $xlNone
$ws.Range("I5").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("I5").Borders.Item($xlEdgeBottom).Weight = $xlThin
$ws.Range("I5").Borders.Item($xlEdgeBottom).ColorIndex = -4105

# ---------------------------------------------------------------------
# 7. Row 6 (was row 5) - Source note, merged A6:H6
# ---------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 27.75
$a6 = $ws.Range("A6")
$a6.Font.Name = "Arial"
$a6.Font.Size = 9
$a6.Font.Bold = $false
$a6.Font.Underline = $false
$a6.Interior.ThemeColor = 0
$a6.Interior.Pattern = 1
$a6.HorizontalAlignment = $xlLeft
$a6.VerticalAlignment = $xlCenter
$a6.WrapText = $true
$a6.Borders.LineStyle = $xlNone

$b6h6 = $ws.Range("B6:H6")
$b6h6.Font.Name = "Arial"
$b6h6.Font.Size = 9
$b6h6.Font.Bold = $false
$b6h6.Interior.ThemeColor = 0
$b6h6.Interior.Pattern = 1
$b6h6.HorizontalAlignment = $xlLeft
$b6h6.VerticalAlignment = $xlCenter
$b6h6.WrapText = $true
$b6h6.Borders.LineStyle = $xlNone
$b6h6.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$b6h6.Borders.Item($xlEdgeTop).Weight = $xlThin
$b6h6.Borders.Item($xlEdgeTop).ColorIndex = -4105

$ws.Range("A6:H6").Merge()

# ---------------------------------------------------------------------
# 8. Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# 9. Selection state to match the saved workbook view
# ---------------------------------------------------------------------
$ws.Range("A1:I1").Select()
